$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3194513333333333
$ws.Range("H2").Value = 0.9583539999999999
$ws.Range("I2").Value = 0.01095865642710367
$ws.Range("J2").Value = 0.01095865642710367
$ws.Range("M2").Value = 0.2901893333333334
$ws.Range("N2").Value = 0.870568
$ws.Range("O2").Value = 0.03429389578125064
$ws.Range("P2").Value = 0.03429389578125064
$ws.Range("Q2").Value = 0.09270136945244445
$ws.Range("R2").Value = 0.8343123250719999
$ws.Range("S2").Value = 0.0003758150214136259
$ws.Range("T2").Value = 0.0003758150214136259

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3194513333333333
$ws.Range("H3").Value = 0.9583539999999999
$ws.Range("I3").Value = 0.01095865642710367
$ws.Range("J3").Value = 0.01095865642710367
$ws.Range("O3").Value = 0.8402845891331153
$ws.Range("P3").Value = 0.8402845891331153
$ws.Range("Q3").Value = 2.271411000934222
$ws.Range("R3").Value = 20.442699008408
$ws.Range("S3").Value = 0.009208390113299783
$ws.Range("T3").Value = 0.009208390113299785

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3194513333333333
$ws.Range("H4").Value = 0.9583539999999999
$ws.Range("I4").Value = 0.01095865642710367
$ws.Range("J4").Value = 0.01095865642710367
$ws.Range("O4").Value = 0.1254215150856341
$ws.Range("P4").Value = 0.1254215150856341
$ws.Range("Q4").Value = 0.3390325287451111
$ws.Range("R4").Value = 3.051292758706
$ws.Range("S4").Value = 0.001374451292390265
$ws.Range("T4").Value = 0.001374451292390264

# Row 5
$ws.Range("I5").Value = 0.9713235907985359
$ws.Range("J5").Value = 0.971323590798536
$ws.Range("M5").Value = 0.2901893333333334
$ws.Range("N5").Value = 0.870568
$ws.Range("O5").Value = 0.03429389578125064
$ws.Range("P5").Value = 0.03429389578125064
$ws.Range("Q5").Value = 8.216611922040888
$ws.Range("R5").Value = 73.94950729836799
$ws.Range("S5").Value = 0.03331046999271513
$ws.Range("T5").Value = 0.03331046999271513

# Row 6
$ws.Range("I6").Value = 0.9713235907985359
$ws.Range("J6").Value = 0.971323590798536
$ws.Range("O6").Value = 0.8402845891331153
$ws.Range("P6").Value = 0.8402845891331153
$ws.Range("S6").Value = 0.81618824440945
$ws.Range("T6").Value = 0.81618824440945

# Row 7
$ws.Range("I7").Value = 0.9713235907985359
$ws.Range("J7").Value = 0.971323590798536
$ws.Range("O7").Value = 0.1254215150856341
$ws.Range("P7").Value = 0.1254215150856341
$ws.Range("Q7").Value = 30.05024342251822
$ws.Range("R7").Value = 270.452190802664
$ws.Range("S7").Value = 0.1218248763963709
$ws.Range("T7").Value = 0.1218248763963709

# Row 8
$ws.Range("I8").Value = 0.01771775277436037
$ws.Range("J8").Value = 0.01771775277436037
$ws.Range("M8").Value = 0.2901893333333334
$ws.Range("N8").Value = 0.870568
$ws.Range("O8").Value = 0.03429389578125064
$ws.Range("P8").Value = 0.03429389578125064
$ws.Range("Q8").Value = 0.149877857448
$ws.Range("R8").Value = 1.348900717032
$ws.Range("S8").Value = 0.000607610767121879
$ws.Range("T8").Value = 0.000607610767121879

# Row 9
$ws.Range("I9").Value = 0.01771775277436037
$ws.Range("J9").Value = 0.01771775277436037
$ws.Range("O9").Value = 0.8402845891331153
$ws.Range("P9").Value = 0.8402845891331153
$ws.Range("S9").Value = 0.01488795461036552
$ws.Range("T9").Value = 0.01488795461036552

# Row 10
$ws.Range("I10").Value = 0.01771775277436037
$ws.Range("J10").Value = 0.01771775277436037
$ws.Range("O10").Value = 0.1254215150856341
$ws.Range("P10").Value = 0.1254215150856341
$ws.Range("Q10").Value = 0.548141514129
$ws.Range("S10").Value = 0.002222187396872975
$ws.Range("T10").Value = 0.002222187396872975

